# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.133.33'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.890.93'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.75'
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5215'
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3756'
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07262'
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.09'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8988'
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08154'
$ws.Range("E12").Value = '  +5.97%  '
$ws.Range("D13").Value = '1.929.72'
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.27'
$ws.Range("E14").Value = '  +1.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.281'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008568'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.56'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '27.160.00'
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.079'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.69'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.397'
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '147.60'
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.285'
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.17'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.737'
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.94'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.782'
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.857'
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09216'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05050'
$ws.Range("E32").Value = '  -0.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7886'
$ws.Range("E33").Value = '  -3.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.213'
$ws.Range("E34").Value = '  -2.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.425'
$ws.Range("E35").Value = '  +3.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.970'
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.580'
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5692'
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01984'
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.011'
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.552'
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.11'
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1516'
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4854'
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.08'
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.621'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.06'
$ws.Range("E49").Value = '  +1.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.38'
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05930'
$ws.Range("E51").Value = '  -0.08%  '
